$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-07-28 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-07-29 Saturday", 2) | Out-Null
$d.Content.Find.Execute("27+42=", $true, $false, $false, $false, $false, $true, 1, $false, "2+65=", 2) | Out-Null
$d.Content.Find.Execute("74-36=", $true, $false, $false, $false, $false, $true, 1, $false, "52+17=", 2) | Out-Null
$d.Content.Find.Execute("73-66=", $true, $false, $false, $false, $false, $true, 1, $false, "43+7=", 2) | Out-Null
$d.Content.Find.Execute("49-48=", $true, $false, $false, $false, $false, $true, 1, $false, "27+13=", 2) | Out-Null
$d.Content.Find.Execute("13+11=", $true, $false, $false, $false, $false, $true, 1, $false, "59-20=", 2) | Out-Null
$d.Content.Find.Execute("2+88=", $true, $false, $false, $false, $false, $true, 1, $false, "64-14=", 2) | Out-Null
$d.Content.Find.Execute("28+28=", $true, $false, $false, $false, $false, $true, 1, $false, "69-36=", 2) | Out-Null
$d.Content.Find.Execute("30+44=", $true, $false, $false, $false, $false, $true, 1, $false, "26-1=", 2) | Out-Null
$d.Content.Find.Execute("52-24=", $true, $false, $false, $false, $false, $true, 1, $false, "49+23=", 2) | Out-Null
$d.Content.Find.Execute("59-14=", $true, $false, $false, $false, $false, $true, 1, $false, "44+52=", 2) | Out-Null
$d.Content.Find.Execute("74-72=", $true, $false, $false, $false, $false, $true, 1, $false, "52+19=", 2) | Out-Null
$d.Content.Find.Execute("39+12=", $true, $false, $false, $false, $false, $true, 1, $false, "85-24=", 2) | Out-Null
$d.Content.Find.Execute("65+16=", $true, $false, $false, $false, $false, $true, 1, $false, "30+28=", 2) | Out-Null
$d.Content.Find.Execute("94-71=", $true, $false, $false, $false, $false, $true, 1, $false, "42-17=", 2) | Out-Null
$d.Content.Find.Execute("91-67=", $true, $false, $false, $false, $false, $true, 1, $false, "78-54=", 2) | Out-Null
$d.Content.Find.Execute("78-41=", $true, $false, $false, $false, $false, $true, 1, $false, "59+17=", 2) | Out-Null
$d.Content.Find.Execute("12+78=", $true, $false, $false, $false, $false, $true, 1, $false, "25-1=", 2) | Out-Null
$d.Content.Find.Execute("54+36=", $true, $false, $false, $false, $false, $true, 1, $false, "76-33=", 2) | Out-Null
$d.Content.Find.Execute("48-24=", $true, $false, $false, $false, $false, $true, 1, $false, "65-35=", 2) | Out-Null
$d.Content.Find.Execute("71-41=", $true, $false, $false, $false, $false, $true, 1, $false, "93-12=", 2) | Out-Null
$d.Content.Find.Execute("31+17=", $true, $false, $false, $false, $false, $true, 1, $false, "61-50=", 2) | Out-Null
$d.Content.Find.Execute("86+1=", $true, $false, $false, $false, $false, $true, 1, $false, "66-55=", 2) | Out-Null
$d.Content.Find.Execute("30+32=", $true, $false, $false, $false, $false, $true, 1, $false, "43+32=", 2) | Out-Null
$d.Content.Find.Execute("26+70=", $true, $false, $false, $false, $false, $true, 1, $false, "12+30=", 2) | Out-Null
$d.Content.Find.Execute("75+19=", $true, $false, $false, $false, $false, $true, 1, $false, "70-13=", 2) | Out-Null
$d.Content.Find.Execute("72-8=", $true, $false, $false, $false, $false, $true, 1, $false, "98-84=", 2) | Out-Null
$d.Content.Find.Execute("58-52=", $true, $false, $false, $false, $false, $true, 1, $false, "9+2=", 2) | Out-Null
$d.Content.Find.Execute("87-11=", $true, $false, $false, $false, $false, $true, 1, $false, "30+46=", 2) | Out-Null
$d.Content.Find.Execute("55-11=", $true, $false, $false, $false, $false, $true, 1, $false, "39-16=", 2) | Out-Null
$d.Content.Find.Execute("65+34=", $true, $false, $false, $false, $false, $true, 1, $false, "10+61=", 2) | Out-Null
$d.Content.Find.Execute("10+47=", $true, $false, $false, $false, $false, $true, 1, $false, "47-41=", 2) | Out-Null
$d.Content.Find.Execute("81-31=", $true, $false, $false, $false, $false, $true, 1, $false, "90-84=", 2) | Out-Null
$d.Content.Find.Execute("69+10=", $true, $false, $false, $false, $false, $true, 1, $false, "68-37=", 2) | Out-Null
$d.Content.Find.Execute("41+23=", $true, $false, $false, $false, $false, $true, 1, $false, "22-13=", 2) | Out-Null
$d.Content.Find.Execute("42+48=", $true, $false, $false, $false, $false, $true, 1, $false, "4+64=", 2) | Out-Null
$d.Content.Find.Execute("42+13=", $true, $false, $false, $false, $false, $true, 1, $false, "10+13=", 2) | Out-Null
$d.Content.Find.Execute("23+75=", $true, $false, $false, $false, $false, $true, 1, $false, "12+18=", 2) | Out-Null
$d.Content.Find.Execute("23+67=", $true, $false, $false, $false, $false, $true, 1, $false, "65-25=", 2) | Out-Null
$d.Content.Find.Execute("41+45=", $true, $false, $false, $false, $false, $true, 1, $false, "71-6=", 2) | Out-Null
$d.Content.Find.Execute("96-60=", $true, $false, $false, $false, $false, $true, 1, $false, "40+22=", 2) | Out-Null
$d.Content.Find.Execute("79+2=", $true, $false, $false, $false, $false, $true, 1, $false, "55-20=", 2) | Out-Null
$d.Content.Find.Execute("19+62=", $true, $false, $false, $false, $false, $true, 1, $false, "20+35=", 2) | Out-Null
$d.Content.Find.Execute("38+42=", $true, $false, $false, $false, $false, $true, 1, $false, "66-39=", 2) | Out-Null
$d.Content.Find.Execute("22+40=", $true, $false, $false, $false, $false, $true, 1, $false, "82-55=", 2) | Out-Null
$d.Content.Find.Execute("94-15=", $true, $false, $false, $false, $false, $true, 1, $false, "8+69=", 2) | Out-Null
$d.Content.Find.Execute("37+10=", $true, $false, $false, $false, $false, $true, 1, $false, "83+13=", 2) | Out-Null
$d.Content.Find.Execute("3+29=", $true, $false, $false, $false, $false, $true, 1, $false, "35-17=", 2) | Out-Null
$d.Content.Find.Execute("89-22=", $true, $false, $false, $false, $false, $true, 1, $false, "59+2=", 2) | Out-Null
$d.Content.Find.Execute("56+23=", $true, $false, $false, $false, $false, $true, 1, $false, "65+13=", 2) | Out-Null
$d.Content.Find.Execute("33+43=", $true, $false, $false, $false, $false, $true, 1, $false, "0+30=", 2) | Out-Null
$d.Content.Find.Execute("56+13=", $true, $false, $false, $false, $false, $true, 1, $false, "97-90=", 2) | Out-Null
$d.Content.Find.Execute("44+36=", $true, $false, $false, $false, $false, $true, 1, $false, "76-17=", 2) | Out-Null
$d.Content.Find.Execute("72-52=", $true, $false, $false, $false, $false, $true, 1, $false, "89-18=", 2) | Out-Null
$d.Content.Find.Execute("29+42=", $true, $false, $false, $false, $false, $true, 1, $false, "69-50=", 2) | Out-Null
$d.Content.Find.Execute("98-87=", $true, $false, $false, $false, $false, $true, 1, $false, "95-1=", 2) | Out-Null
$d.Content.Find.Execute("83-41=", $true, $false, $false, $false, $false, $true, 1, $false, "21-10=", 2) | Out-Null
$d.Content.Find.Execute("46+53=", $true, $false, $false, $false, $false, $true, 1, $false, "26+1=", 2) | Out-Null
$d.Content.Find.Execute("69-15=", $true, $false, $false, $false, $false, $true, 1, $false, "3+94=", 2) | Out-Null
$d.Content.Find.Execute("4+5=", $true, $false, $false, $false, $false, $true, 1, $false, "61-21=", 2) | Out-Null
$d.Content.Find.Execute("42-40=", $true, $false, $false, $false, $false, $true, 1, $false, "85-83=", 2) | Out-Null
$d.Content.Find.Execute("86-76=", $true, $false, $false, $false, $false, $true, 1, $false, "39+31=", 2) | Out-Null
$d.Content.Find.Execute("77+21=", $true, $false, $false, $false, $false, $true, 1, $false, "19+2=", 2) | Out-Null
$d.Content.Find.Execute("34+18=", $true, $false, $false, $false, $false, $true, 1, $false, "76-43=", 2) | Out-Null
$d.Content.Find.Execute("65-55=", $true, $false, $false, $false, $false, $true, 1, $false, "36+24=", 2) | Out-Null
$d.Content.Find.Execute("54-54=", $true, $false, $false, $false, $false, $true, 1, $false, "24-17=", 2) | Out-Null
$d.Content.Find.Execute("92-89=", $true, $false, $false, $false, $false, $true, 1, $false, "56+6=", 2) | Out-Null
$d.Content.Find.Execute("8+1=", $true, $false, $false, $false, $false, $true, 1, $false, "12+12=", 2) | Out-Null
$d.Content.Find.Execute("67-16=", $true, $false, $false, $false, $false, $true, 1, $false, "67-47=", 2) | Out-Null
$d.Content.Find.Execute("6+11=", $true, $false, $false, $false, $false, $true, 1, $false, "85-19=", 2) | Out-Null
$d.Content.Find.Execute("30-24=", $true, $false, $false, $false, $false, $true, 1, $false, "66-60=", 2) | Out-Null
$d.Content.Find.Execute("20+56=", $true, $false, $false, $false, $false, $true, 1, $false, "24+45=", 2) | Out-Null
$d.Content.Find.Execute("39+15=", $true, $false, $false, $false, $false, $true, 1, $false, "43+30=", 2) | Out-Null
$d.Content.Find.Execute("33+32=", $true, $false, $false, $false, $false, $true, 1, $false, "17+55=", 2) | Out-Null
$d.Content.Find.Execute("59+19=", $true, $false, $false, $false, $false, $true, 1, $false, "66-27=", 2) | Out-Null
$d.Content.Find.Execute("74-71=", $true, $false, $false, $false, $false, $true, 1, $false, "21+58=", 2) | Out-Null
$d.Content.Find.Execute("96-17=", $true, $false, $false, $false, $false, $true, 1, $false, "94-54=", 2) | Out-Null
$d.Content.Find.Execute("15-6=", $true, $false, $false, $false, $false, $true, 1, $false, "72+11=", 2) | Out-Null
$d.Content.Find.Execute("89-56=", $true, $false, $false, $false, $false, $true, 1, $false, "98-1=", 2) | Out-Null
$d.Content.Find.Execute("2+92=", $true, $false, $false, $false, $false, $true, 1, $false, "12+33=", 2) | Out-Null
$d.Content.Find.Execute("2+81=", $true, $false, $false, $false, $false, $true, 1, $false, "72-31=", 2) | Out-Null
$d.Content.Find.Execute("39-26=", $true, $false, $false, $false, $false, $true, 1, $false, "42-19=", 2) | Out-Null
$d.Content.Find.Execute("34+25=", $true, $false, $false, $false, $false, $true, 1, $false, "21-0=", 2) | Out-Null
$d.Content.Find.Execute("98-32=", $true, $false, $false, $false, $false, $true, 1, $false, "50+25=", 2) | Out-Null
$d.Content.Find.Execute("18+26=", $true, $false, $false, $false, $false, $true, 1, $false, "85-42=", 2) | Out-Null
$d.Content.Find.Execute("33+63=", $true, $false, $false, $false, $false, $true, 1, $false, "99-30=", 2) | Out-Null
$d.Content.Find.Execute("54+30=", $true, $false, $false, $false, $false, $true, 1, $false, "64+13=", 2) | Out-Null
$d.Content.Find.Execute("66-6=", $true, $false, $false, $false, $false, $true, 1, $false, "81-38=", 2) | Out-Null
$d.Content.Find.Execute("19+43=", $true, $false, $false, $false, $false, $true, 1, $false, "18-16=", 2) | Out-Null
$d.Content.Find.Execute("10-7=", $true, $false, $false, $false, $false, $true, 1, $false, "53+30=", 2) | Out-Null
$d.Content.Find.Execute("95-74=", $true, $false, $false, $false, $false, $true, 1, $false, "39+32=", 2) | Out-Null
$d.Content.Find.Execute("10+46=", $true, $false, $false, $false, $false, $true, 1, $false, "13+11=", 2) | Out-Null
$d.Content.Find.Execute("51-46=", $true, $false, $false, $false, $false, $true, 1, $false, "7+52=", 2) | Out-Null
$d.Content.Find.Execute("20+38=", $true, $false, $false, $false, $false, $true, 1, $false, "79-44=", 2) | Out-Null
$d.Content.Find.Execute("19+65=", $true, $false, $false, $false, $false, $true, 1, $false, "66+8=", 2) | Out-Null
$d.Content.Find.Execute("92-23=", $true, $false, $false, $false, $false, $true, 1, $false, "91-61=", 2) | Out-Null
$d.Content.Find.Execute("7+15=", $true, $false, $false, $false, $false, $true, 1, $false, "77-63=", 2) | Out-Null
$d.Content.Find.Execute("89-42=", $true, $false, $false, $false, $false, $true, 1, $false, "88-64=", 2) | Out-Null
$d.Content.Find.Execute("6+5=", $true, $false, $false, $false, $false, $true, 1, $false, "62+24=", 2) | Out-Null
$d.Content.Find.Execute("8+56=", $true, $false, $false, $false, $false, $true, 1, $false, "63-52=", 2) | Out-Null
$d.Content.Find.Execute("18+7=", $true, $false, $false, $false, $false, $true, 1, $false, "72+16=", 2) | Out-Null
